$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells Q1/R1: copy style from P1 (same header formatting) before setting values
$ws.Range("P1").Copy($ws.Range("Q1:R1"))

$ws.Range("B1").Value = 'Customer type'
$ws.Range("C1").Value = 'Load level average [MW]'
$ws.Range("D1").Value = 'Load point peak [MW]'
$ws.Range("E1").Value = 'Number of customers'
$ws.Range("F1").Value = 'U'
$ws.Range("G1").Value = 'nrOfFaults'
$ws.Range("H1").Value = 'R'
$ws.Range("I1").Value = 'Lambda'
$ws.Range("J1").Value = 'SAIFI'
$ws.Range("K1").Value = 'SAIDI'
$ws.Range("L1").Value = 'CAIDI'
$ws.Range("M1").Value = 'EENS'
$ws.Range("N1").Value = 'nr of simulations'
$ws.Range("O1").Value = 'provided beta'
$ws.Range("P1").Value = 'calculated beta'
$ws.Range("Q1").Value = 'EENS 95% CI'
$ws.Range("R1").Value = 'EENS 99% CI'
$ws.Range("A2").Value = 'LP1'
$ws.Range("B2").Value = 'residential'
$ws.Range("C2").Value = 0.545
$ws.Range("D2").Value = 0.8869
$ws.Range("E2").Value = 220
$ws.Range("F2").Value = 2.288799477429029
$ws.Range("G2").Value = 16167
$ws.Range("H2").Value = 3.111193005256407
$ws.Range("I2").Value = 0.7356661812886786
$ws.Range("J2").Value = 161.8465598835093
$ws.Range("K2").Value = 503.5358850343864
$ws.Range("L2").Value = 684.4624611564096
$ws.Range("M2").Value = 1.247395715198821
$ws.Range("A3").Value = 'LP2'
$ws.Range("B3").Value = 'residential'
$ws.Range("C3").Value = 0.545
$ws.Range("D3").Value = 0.8869
$ws.Range("E3").Value = 220
$ws.Range("F3").Value = 2.349138018721725
$ws.Range("G3").Value = 16399
$ws.Range("H3").Value = 3.148036898556536
$ws.Range("I3").Value = 0.7462231525300328
$ws.Range("J3").Value = 164.1690935566072
$ws.Range("K3").Value = 516.8103641187796
$ws.Range("L3").Value = 692.5681176824379
$ws.Range("M3").Value = 1.280280220203341
$ws.Range("A4").Value = 'LP3'
$ws.Range("B4").Value = 'residential'
$ws.Range("C4").Value = 0.545
$ws.Range("D4").Value = 0.8869
$ws.Range("E4").Value = 220
$ws.Range("F4").Value = 2.303132337188917
$ws.Range("G4").Value = 16173
$ws.Range("H4").Value = 3.129514390778684
$ws.Range("I4").Value = 0.7359392064069894
$ws.Range("J4").Value = 161.9066254095377
$ws.Range("K4").Value = 506.6891141815618
$ws.Range("L4").Value = 688.4931659713104
$ws.Range("M4").Value = 1.25520712376796
$ws.Range("A5").Value = 'LP4'
$ws.Range("B5").Value = 'residential'
$ws.Range("C5").Value = 0.545
$ws.Range("D5").Value = 0.8869
$ws.Range("E5").Value = 220
$ws.Range("F5").Value = 2.312062703893669
$ws.Range("G5").Value = 16406
$ws.Range("H5").Value = 3.097030963109062
$ws.Range("I5").Value = 0.7465416818347288
$ws.Range("J5").Value = 164.2391700036403
$ws.Range("K5").Value = 508.6537948566071
$ws.Range("L5").Value = 681.3468118839935
$ws.Range("M5").Value = 1.26007417362205
$ws.Range("A6").Value = 'LP5'
$ws.Range("B6").Value = 'residential'
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.8137
$ws.Range("E6").Value = 200
$ws.Range("F6").Value = 2.334387154306878
$ws.Range("G6").Value = 16344
$ws.Range("H6").Value = 3.138796628918744
$ws.Range("I6").Value = 0.7437204222788496
$ws.Range("J6").Value = 148.7440844557699
$ws.Range("K6").Value = 466.8774308613756
$ws.Range("L6").Value = 627.7593257837487
$ws.Range("M6").Value = 1.167193577153439
$ws.Range("A7").Value = 'LP6'
$ws.Range("B7").Value = 'commercial'
$ws.Range("C7").Value = 0.415
$ws.Range("D7").Value = 0.6714
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 2.364057890998176
$ws.Range("G7").Value = 16507
$ws.Range("H7").Value = 3.147303338739682
$ws.Range("I7").Value = 0.7511376046596286
$ws.Range("J7").Value = 7.511376046596286
$ws.Range("K7").Value = 23.64057890998176
$ws.Range("L7").Value = 31.47303338739682
$ws.Range("M7").Value = 0.9810840247642432
$ws.Range("A8").Value = 'LP7'
$ws.Range("B8").Value = 'commercial'
$ws.Range("C8").Value = 0.415
$ws.Range("D8").Value = 0.6714
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 2.308495903996086
$ws.Range("G8").Value = 16366
$ws.Range("H8").Value = 3.0998109486874
$ws.Range("I8").Value = 0.7447215143793229
$ws.Range("J8").Value = 7.447215143793229
$ws.Range("K8").Value = 23.08495903996086
$ws.Range("L8").Value = 30.998109486874
$ws.Range("M8").Value = 0.9580258001583757
$ws.Range("A9").Value = 'LP8'
$ws.Range("B9").Value = 'small user'
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1.63
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 2.13779155396402
$ws.Range("G9").Value = 15824
$ws.Range("H9").Value = 2.968914761748819
$ws.Range("I9").Value = 0.720058245358573
$ws.Range("J9").Value = 0.720058245358573
$ws.Range("K9").Value = 2.13779155396402
$ws.Range("L9").Value = 2.968914761748819
$ws.Range("M9").Value = 2.13779155396402
$ws.Range("A10").Value = 'LP9'
$ws.Range("B10").Value = 'small user'
$ws.Range("C10").Value = 1.5
$ws.Range("D10").Value = 2.445
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2.194773954555247
$ws.Range("G10").Value = 16102
$ws.Range("H10").Value = 2.995426184654459
$ws.Range("I10").Value = 0.732708409173644
$ws.Range("J10").Value = 0.732708409173644
$ws.Range("K10").Value = 2.194773954555247
$ws.Range("L10").Value = 2.995426184654459
$ws.Range("M10").Value = 3.29216093183287
$ws.Range("A11").Value = 'LP10'
$ws.Range("B11").Value = 'small user'
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1.63
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 2.223771496129151
$ws.Range("G11").Value = 16188
$ws.Range("H11").Value = 3.018878329561046
$ws.Range("I11").Value = 0.7366217692027667
$ws.Range("J11").Value = 0.7366217692027667
$ws.Range("K11").Value = 2.223771496129151
$ws.Range("L11").Value = 3.018878329561046
$ws.Range("M11").Value = 2.223771496129151
$ws.Range("A12").Value = 'LP11'
$ws.Range("B12").Value = 'residential'
$ws.Range("C12").Value = 0.545
$ws.Range("D12").Value = 0.8869
$ws.Range("E12").Value = 220
$ws.Range("F12").Value = 2.325472194264243
$ws.Range("G12").Value = 16376
$ws.Range("H12").Value = 3.120699617803554
$ws.Range("I12").Value = 0.7451765562431744
$ws.Range("J12").Value = 163.9388423734983
$ws.Range("K12").Value = 511.6038827381334
$ws.Range("L12").Value = 686.5539159167819
$ws.Range("M12").Value = 1.267382345874013
$ws.Range("A13").Value = 'LP12'
$ws.Range("B13").Value = 'residential'
$ws.Range("C13").Value = 0.545
$ws.Range("D13").Value = 0.8869
$ws.Range("E13").Value = 220
$ws.Range("F13").Value = 2.338461779785088
$ws.Range("G13").Value = 16426
$ws.Range("H13").Value = 3.128578842844095
$ws.Range("I13").Value = 0.7474517655624318
$ws.Range("J13").Value = 164.439388423735
$ws.Range("K13").Value = 514.4615915527194
$ws.Range("L13").Value = 688.2873454257008
$ws.Range("M13").Value = 1.274461669982873
$ws.Range("A14").Value = 'LP13'
$ws.Range("B14").Value = 'residential'
$ws.Range("C14").Value = 0.545
$ws.Range("D14").Value = 0.8869
$ws.Range("E14").Value = 220
$ws.Range("F14").Value = 2.326345192294136
$ws.Range("G14").Value = 16373
$ws.Range("H14").Value = 3.122443165324372
$ws.Range("I14").Value = 0.7450400436840189
$ws.Range("J14").Value = 163.9088096104842
$ws.Range("K14").Value = 511.79594230471
$ws.Range("L14").Value = 686.9374963713618
$ws.Range("M14").Value = 1.267858129800304
$ws.Range("A15").Value = 'LP14'
$ws.Range("B15").Value = 'residential'
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 0.8137
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 2.300497468537783
$ws.Range("G15").Value = 16206
$ws.Range("H15").Value = 3.119568824422209
$ws.Range("I15").Value = 0.7374408445576993
$ws.Range("J15").Value = 147.4881689115399
$ws.Range("K15").Value = 460.0994937075566
$ws.Range("L15").Value = 623.9137648844418
$ws.Range("M15").Value = 1.150248734268891
$ws.Range("A16").Value = 'LP15'
$ws.Range("B16").Value = 'residential'
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 0.8137
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 2.334813813171912
$ws.Range("G16").Value = 16366
$ws.Range("H16").Value = 3.135150211307951
$ws.Range("I16").Value = 0.7447215143793229
$ws.Range("J16").Value = 148.9443028758646
$ws.Range("K16").Value = 466.9627626343823
$ws.Range("L16").Value = 627.0300422615903
$ws.Range("M16").Value = 1.167406906585956
$ws.Range("A17").Value = 'LP16'
$ws.Range("B17").Value = 'commercial'
$ws.Range("C17").Value = 0.415
$ws.Range("D17").Value = 0.6714
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 2.276229104842548
$ws.Range("G17").Value = 16176
$ws.Range("H17").Value = 3.092384446588764
$ws.Range("I17").Value = 0.7360757189661449
$ws.Range("J17").Value = 7.360757189661449
$ws.Range("K17").Value = 22.76229104842548
$ws.Range("L17").Value = 30.92384446588764
$ws.Range("M17").Value = 0.9446350785096574
$ws.Range("A18").Value = 'LP17'
$ws.Range("B18").Value = 'commercial'
$ws.Range("C18").Value = 0.415
$ws.Range("D18").Value = 0.6714
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 2.353188777367233
$ws.Range("G18").Value = 16453
$ws.Range("H18").Value = 3.143115332852507
$ws.Range("I18").Value = 0.7486803785948307
$ws.Range("J18").Value = 7.486803785948307
$ws.Range("K18").Value = 23.53188777367233
$ws.Range("L18").Value = 31.43115332852507
$ws.Range("M18").Value = 0.9765733426074015
$ws.Range("A19").Value = 'LP18'
$ws.Range("B19").Value = 'residential'
$ws.Range("C19").Value = 0.545
$ws.Range("D19").Value = 0.8869
$ws.Range("E19").Value = 220
$ws.Range("F19").Value = 2.123146765109679
$ws.Range("G19").Value = 11747
$ws.Range("H19").Value = 3.97193098749045
$ws.Range("I19").Value = 0.5345376774663269
$ws.Range("J19").Value = 117.5982890425919
$ws.Range("K19").Value = 467.0922883241295
$ws.Range("L19").Value = 873.8248172478991
$ws.Range("M19").Value = 1.157114986984775
$ws.Range("A20").Value = 'LP19'
$ws.Range("B20").Value = 'residential'
$ws.Range("C20").Value = 0.545
$ws.Range("D20").Value = 0.8869
$ws.Range("E20").Value = 220
$ws.Range("F20").Value = 2.061480692061452
$ws.Range("G20").Value = 11439
$ws.Range("H20").Value = 3.960407351057126
$ws.Range("I20").Value = 0.5205223880597015
$ws.Range("J20").Value = 114.5149253731343
$ws.Range("K20").Value = 453.5257522535194
$ws.Range("L20").Value = 871.2896172325676
$ws.Range("M20").Value = 1.123506977173491
$ws.Range("A21").Value = 'LP20'
$ws.Range("B21").Value = 'residential'
$ws.Range("C21").Value = 0.545
$ws.Range("D21").Value = 0.8869
$ws.Range("E21").Value = 220
$ws.Range("F21").Value = 2.103584627195791
$ws.Range("G21").Value = 11681
$ws.Range("H21").Value = 3.957570051130443
$ws.Range("I21").Value = 0.5315344011649071
$ws.Range("J21").Value = 116.9375682562796
$ws.Range("K21").Value = 462.788617983074
$ws.Range("L21").Value = 870.6654112486974
$ws.Range("M21").Value = 1.146453621821706
$ws.Range("A22").Value = 'LP21'
$ws.Range("B22").Value = 'residential'
$ws.Range("C22").Value = 0.545
$ws.Range("D22").Value = 0.8869
$ws.Range("E22").Value = 220
$ws.Range("F22").Value = 2.130183722536171
$ws.Range("G22").Value = 11706
$ws.Range("H22").Value = 3.999053262126679
$ws.Range("I22").Value = 0.5326720058245359
$ws.Range("J22").Value = 117.1878412813979
$ws.Range("K22").Value = 468.6404189579577
$ws.Range("L22").Value = 879.7917176678694
$ws.Range("M22").Value = 1.160950128782213
$ws.Range("A23").Value = 'LP22'
$ws.Range("B23").Value = 'residential'
$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = 0.8137
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 2.069354984638551
$ws.Range("G23").Value = 11497
$ws.Range("H23").Value = 3.955479267845246
$ws.Range("I23").Value = 0.5231616308700401
$ws.Range("J23").Value = 104.632326174008
$ws.Range("K23").Value = 413.8709969277102
$ws.Range("L23").Value = 791.0958535690492
$ws.Range("M23").Value = 1.034677492319275
$ws.Range("A24").Value = 'LP23'
$ws.Range("B24").Value = 'residential'
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 0.8137
$ws.Range("E24").Value = 200
$ws.Range("F24").Value = 2.114473727487074
$ws.Range("G24").Value = 11720
$ws.Range("H24").Value = 3.964818654885319
$ws.Range("I24").Value = 0.5333090644339279
$ws.Range("J24").Value = 106.6618128867856
$ws.Range("K24").Value = 422.8947454974148
$ws.Range("L24").Value = 792.9637309770638
$ws.Range("M24").Value = 1.057236863743537
$ws.Range("A25").Value = 'LP24'
$ws.Range("B25").Value = 'commercial'
$ws.Range("C25").Value = 0.415
$ws.Range("D25").Value = 0.6714
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 2.154559449981684
$ws.Range("G25").Value = 11742
$ws.Range("H25").Value = 4.032413428104028
$ws.Range("I25").Value = 0.5343101565344012
$ws.Range("J25").Value = 5.343101565344011
$ws.Range("K25").Value = 21.54559449981684
$ws.Range("L25").Value = 40.32413428104028
$ws.Range("M25").Value = 0.894142171742399
$ws.Range("A26").Value = 'LP25'
$ws.Range("B26").Value = 'commercial'
$ws.Range("C26").Value = 0.415
$ws.Range("D26").Value = 0.6714
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 2.07045324812359
$ws.Range("G26").Value = 11502
$ws.Range("H26").Value = 3.955858162125197
$ws.Range("I26").Value = 0.5233891518019658
$ws.Range("J26").Value = 5.233891518019657
$ws.Range("K26").Value = 20.7045324812359
$ws.Range("L26").Value = 39.55858162125197
$ws.Range("M26").Value = 0.8592380979712899
$ws.Range("A27").Value = 'LP26'
$ws.Range("B27").Value = 'small user'
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1.63
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1.981760394344859
$ws.Range("G27").Value = 11344
$ws.Range("H27").Value = 3.839136673670895
$ws.Range("I27").Value = 0.5161994903531125
$ws.Range("J27").Value = 0.5161994903531125
$ws.Range("K27").Value = 1.981760394344859
$ws.Range("L27").Value = 3.839136673670895
$ws.Range("M27").Value = 1.981760394344859
$ws.Range("A28").Value = 'LP27'
$ws.Range("B28").Value = 'small user'
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1.63
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 2.004778934152148
$ws.Range("G28").Value = 11492
$ws.Range("H28").Value = 3.833712309165299
$ws.Range("I28").Value = 0.5229341099381143
$ws.Range("J28").Value = 0.5229341099381143
$ws.Range("K28").Value = 2.004778934152148
$ws.Range("L28").Value = 3.833712309165299
$ws.Range("M28").Value = 2.004778934152148
$ws.Range("A29").Value = 'LP28'
$ws.Range("B29").Value = 'small user'
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1.63
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 1.940447001897347
$ws.Range("G29").Value = 11179
$ws.Range("H29").Value = 3.814586574263895
$ws.Range("I29").Value = 0.5086912995995632
$ws.Range("J29").Value = 0.5086912995995632
$ws.Range("K29").Value = 1.940447001897347
$ws.Range("L29").Value = 3.814586574263895
$ws.Range("M29").Value = 1.940447001897347
$ws.Range("A30").Value = 'LP29'
$ws.Range("B30").Value = 'small user'
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1.63
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 1.939473328760741
$ws.Range("G30").Value = 11160
$ws.Range("H30").Value = 3.819163608677961
$ws.Range("I30").Value = 0.5078267200582454
$ws.Range("J30").Value = 0.5078267200582454
$ws.Range("K30").Value = 1.939473328760741
$ws.Range("L30").Value = 3.819163608677961
$ws.Range("M30").Value = 1.939473328760741
$ws.Range("A31").Value = 'LP30'
$ws.Range("B31").Value = 'small user'
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 1.63
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 1.989990560565909
$ws.Range("G31").Value = 11406
$ws.Range("H31").Value = 3.834125246273576
$ws.Range("I31").Value = 0.5190207499089916
$ws.Range("J31").Value = 0.5190207499089916
$ws.Range("K31").Value = 1.989990560565909
$ws.Range("L31").Value = 3.834125246273576
$ws.Range("M31").Value = 1.989990560565909
$ws.Range("A32").Value = 'LP31'
$ws.Range("B32").Value = 'small user'
$ws.Range("C32").Value = 1.5
$ws.Range("D32").Value = 2.445
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 1.933678469572113
$ws.Range("G32").Value = 11163
$ws.Range("H32").Value = 3.806729198899646
$ws.Range("I32").Value = 0.5079632326174008
$ws.Range("J32").Value = 0.5079632326174008
$ws.Range("K32").Value = 1.933678469572113
$ws.Range("L32").Value = 3.806729198899646
$ws.Range("M32").Value = 2.900517704358169
$ws.Range("A33").Value = 'LP32'
$ws.Range("B33").Value = 'residential'
$ws.Range("C33").Value = 0.545
$ws.Range("D33").Value = 0.8869
$ws.Range("E33").Value = 220
$ws.Range("F33").Value = 2.138917983530538
$ws.Range("G33").Value = 11761
$ws.Range("H33").Value = 3.996672188255003
$ws.Range("I33").Value = 0.535174736075719
$ws.Range("J33").Value = 117.7384419366582
$ws.Range("K33").Value = 470.5619563767183
$ws.Range("L33").Value = 879.2678814161007
$ws.Range("M33").Value = 1.165710301024143
$ws.Range("A34").Value = 'LP33'
$ws.Range("B34").Value = 'residential'
$ws.Range("C34").Value = 0.545
$ws.Range("D34").Value = 0.8869
$ws.Range("E34").Value = 220
$ws.Range("F34").Value = 2.13297544800255
$ws.Range("G34").Value = 11730
$ws.Range("H34").Value = 3.996101316735213
$ws.Range("I34").Value = 0.5337641062977794
$ws.Range("J34").Value = 117.4281033855115
$ws.Range("K34").Value = 469.2545985605611
$ws.Range("L34").Value = 879.1422896817468
$ws.Range("M34").Value = 1.16247161916139
$ws.Range("A35").Value = 'LP34'
$ws.Range("B35").Value = 'residential'
$ws.Range("C35").Value = 0.545
$ws.Range("D35").Value = 0.8869
$ws.Range("E35").Value = 220
$ws.Range("F35").Value = 2.076326848683585
$ws.Range("G35").Value = 11486
$ws.Range("H35").Value = 3.972606549422816
$ws.Range("I35").Value = 0.5226610848198034
$ws.Range("J35").Value = 114.9854386603568
$ws.Range("K35").Value = 456.7919067103887
$ws.Range("L35").Value = 873.9734408730195
$ws.Range("M35").Value = 1.131598132532554
$ws.Range("A36").Value = 'LP35'
$ws.Range("B36").Value = 'residential'
$ws.Range("C36").Value = 0.545
$ws.Range("D36").Value = 0.8869
$ws.Range("E36").Value = 220
$ws.Range("F36").Value = 2.12669777951851
$ws.Range("G36").Value = 11719
$ws.Range("H36").Value = 3.988080075322023
$ws.Range("I36").Value = 0.5332635602475427
$ws.Range("J36").Value = 117.3179832544594
$ws.Range("K36").Value = 467.8735114940723
$ws.Range("L36").Value = 877.377616570845
$ws.Range("M36").Value = 1.159050289837588
$ws.Range("A37").Value = 'LP36'
$ws.Range("B37").Value = 'residential'
$ws.Range("C37").Value = 0.5
$ws.Range("D37").Value = 0.8137
$ws.Range("E37").Value = 200
$ws.Range("F37").Value = 2.092978466538161
$ws.Range("G37").Value = 11539
$ws.Range("H37").Value = 3.986072864255362
$ws.Range("I37").Value = 0.5250728066982162
$ws.Range("J37").Value = 105.0145613396432
$ws.Range("K37").Value = 418.5956933076322
$ws.Range("L37").Value = 797.2145728510725
$ws.Range("M37").Value = 1.04648923326908
$ws.Range("A38").Value = 'LP37'
$ws.Range("B38").Value = 'residential'
$ws.Range("C38").Value = 0.5
$ws.Range("D38").Value = 0.8137
$ws.Range("E38").Value = 200
$ws.Range("F38").Value = 2.148208800140189
$ws.Range("G38").Value = 11779
$ws.Range("H38").Value = 4.007898513615824
$ws.Range("I38").Value = 0.5359938114306516
$ws.Range("J38").Value = 107.1987622861303
$ws.Range("K38").Value = 429.6417600280378
$ws.Range("L38").Value = 801.5797027231649
$ws.Range("M38").Value = 1.074104400070095
$ws.Range("A39").Value = 'LP38'
$ws.Range("B39").Value = 'commercial'
$ws.Range("C39").Value = 0.415
$ws.Range("D39").Value = 0.6714
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 2.088633162806648
$ws.Range("G39").Value = 11481
$ws.Range("H39").Value = 3.997892377479217
$ws.Range("I39").Value = 0.5224335638878776
$ws.Range("J39").Value = 5.224335638878776
$ws.Range("K39").Value = 20.88633162806648
$ws.Range("L39").Value = 39.97892377479217
$ws.Range("M39").Value = 0.8667827625647588
$ws.Range("A40").Value = 'TOTAL'
$ws.Range("C40").Value = 24.58000000000001
$ws.Range("D40").Value = 39.99919999999999
$ws.Range("E40").Value = 4779
$ws.Range("J40").Value = 0.6272694296496327
$ws.Range("K40").Value = 2.206219951765333
$ws.Range("L40").Value = 3.517180732046257
$ws.Range("M40").Value = 52.64804582750083
$ws.Range("N40").Value = 21976
$ws.Range("O40").Value = 0.02
$ws.Range("P40").Value = 0.02084032945293676
$ws.Range("Q40").Value = '(50.4987964686871, 54.799938720014886)'
$ws.Range("R40").Value = '(49.82290268633559, 55.475832502366394)'
